$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '310.28'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-2.25%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '53.72'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '12.56%'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-3.30%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07813'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-1.60%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '4.509'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-1.91%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.360'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '3.31%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.587'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-3.14%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1229'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-3.91%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.2005'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '3.49%'
$ws.Range("B11").Value = 'BitrueCoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.04727'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '1.85%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09408'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.50%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.1046'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.05%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001257'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-4.79%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005805'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-1.07%'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '2,010.15%'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '0.08%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.420'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.49%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3447'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.52%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.989'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-1.08%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1364'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-2.23%'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-0.47%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04176'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.08%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001259'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-4.73%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.003976'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-5.24%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02602'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '-2.07%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05953'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '3.41%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01098'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '1.85%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007913'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.34%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1423'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-0.71%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.008225'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '6.80%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008471'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-0.22%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3128'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-1.05%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00007319'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '5.75%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.38%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05654'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '3.00%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002620'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-34.65%'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.38%'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.38%'
